$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.136.80"
$ws.Range("D3").Value = "2.643.71"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "596.90"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "156.65"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "0.352"
$ws.Range("D13").Value = "28.04"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "3.124.76"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "68.064.22"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "2.643.24"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "11.35"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "362.59"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D21").Value = "4.42"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").Value = "75.08"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "9.70"
$ws.Range("E26").Value = "  -3.81%  "
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "555.38"
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "159.86"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "40.42"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "158.91"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "22.06"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "0.0785"
$ws.Range("E51").Value = "  +0.19%  "
